$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 1.779961347579956
$ws.Cells.Item(2, 5).Value = 2084.208947092691
$ws.Cells.Item(2, 6).Value = 0.09410775573023475
$ws.Cells.Item(2, 7).Value = 0.07685615811998533
$ws.Cells.Item(2, 8).Value = 0.0650379480944331
$ws.Cells.Item(2, 9).Value = 0.06016307334290897
$ws.Cells.Item(2, 10).Value = 0.05398993049699734
$ws.Cells.Item(2, 11).Value = 0.05092966470659509
$ws.Cells.Item(2, 12).Value = 0.04932695777194288
$ws.Cells.Item(2, 13).Value = 0.04732895414439074
$ws.Cells.Item(2, 14).Value = 0.04542350901325164
$ws.Cells.Item(2, 15).Value = 0.0447847250094331
$ws.Cells.Item(2, 16).Value = 0.04374833676317175
$ws.Cells.Item(2, 17).Value = 0.04281507856481386
$ws.Cells.Item(2, 18).Value = 0.0423257174870185
$ws.Cells.Item(2, 19).Value = 0.04184543948415295
$ws.Cells.Item(2, 20).Value = 0.04152481837247183
$ws.Cells.Item(2, 21).Value = 0.04122362290224218
$ws.Cells.Item(2, 22).Value = 0.04107781108793178
$ws.Cells.Item(2, 23).Value = 0.04088138072157833
$ws.Cells.Item(2, 24).Value = 0.04067520002878714
$ws.Cells.Item(2, 25).Value = 0.04062785471915577
$ws.Cells.Item(3, 3).Value = 1.498037338256836
$ws.Cells.Item(3, 5).Value = 2076.194138023489
$ws.Cells.Item(3, 6).Value = 0.0946640925639493
$ws.Cells.Item(3, 7).Value = 0.07547429697127506
$ws.Cells.Item(3, 8).Value = 0.06523312113804133
$ws.Cells.Item(3, 9).Value = 0.05897679766815391
$ws.Cells.Item(3, 10).Value = 0.05203247794618173
$ws.Cells.Item(3, 11).Value = 0.05203247794618173
$ws.Cells.Item(3, 12).Value = 0.04987208825711924
$ws.Cells.Item(3, 13).Value = 0.04684346996895047
$ws.Cells.Item(3, 14).Value = 0.0451714968043502
$ws.Cells.Item(3, 15).Value = 0.0444109924468123
$ws.Cells.Item(3, 16).Value = 0.04370612949872715
$ws.Cells.Item(3, 17).Value = 0.04302194869373173
$ws.Cells.Item(3, 18).Value = 0.04268021211882009
$ws.Cells.Item(3, 19).Value = 0.04203050931441327
$ws.Cells.Item(3, 20).Value = 0.04151662474600519
$ws.Cells.Item(3, 21).Value = 0.04119797139275679
$ws.Cells.Item(3, 22).Value = 0.04076095041551464
$ws.Cells.Item(3, 23).Value = 0.04069124437802291
$ws.Cells.Item(3, 24).Value = 0.04057654674918851
$ws.Cells.Item(3, 25).Value = 0.04047162062423955
$ws.Cells.Item(4, 3).Value = 1.616026163101196
$ws.Cells.Item(4, 5).Value = 2056.398111952029
$ws.Cells.Item(4, 6).Value = 0.09332948515268631
$ws.Cells.Item(4, 7).Value = 0.07723116330518767
$ws.Cells.Item(4, 8).Value = 0.06628219646908315
$ws.Cells.Item(4, 9).Value = 0.05962900013117365
$ws.Cells.Item(4, 10).Value = 0.05506320160501134
$ws.Cells.Item(4, 11).Value = 0.05229243294305164
$ws.Cells.Item(4, 12).Value = 0.04974306139253499
$ws.Cells.Item(4, 13).Value = 0.04757399622798197
$ws.Cells.Item(4, 14).Value = 0.04524831247848512
$ws.Cells.Item(4, 15).Value = 0.04443536272501646
$ws.Cells.Item(4, 16).Value = 0.04324437233606658
$ws.Cells.Item(4, 17).Value = 0.04190995944332584
$ws.Cells.Item(4, 18).Value = 0.04161653922613266
$ws.Cells.Item(4, 19).Value = 0.04092192235605797
$ws.Cells.Item(4, 20).Value = 0.0408366625400271
$ws.Cells.Item(4, 21).Value = 0.04061484300975535
$ws.Cells.Item(4, 22).Value = 0.04046030001575787
$ws.Cells.Item(4, 23).Value = 0.04025115565116848
$ws.Cells.Item(4, 24).Value = 0.04015323457082132
$ws.Cells.Item(4, 25).Value = 0.04008573317645279
$ws.Cells.Item(5, 3).Value = 1.580986261367798
$ws.Cells.Item(5, 5).Value = 2089.528007068355
$ws.Cells.Item(5, 6).Value = 0.09407049365186053
$ws.Cells.Item(5, 7).Value = 0.07366930900259971
$ws.Cells.Item(5, 8).Value = 0.06602704761586128
$ws.Cells.Item(5, 9).Value = 0.05998068221838463
$ws.Cells.Item(5, 10).Value = 0.05640466573466887
$ws.Cells.Item(5, 11).Value = 0.05332802170903994
$ws.Cells.Item(5, 12).Value = 0.0506796524571128
$ws.Cells.Item(5, 13).Value = 0.04824055022335234
$ws.Cells.Item(5, 14).Value = 0.04637531498951879
$ws.Cells.Item(5, 15).Value = 0.04470779480651561
$ws.Cells.Item(5, 16).Value = 0.04402883712844365
$ws.Cells.Item(5, 17).Value = 0.04309453677487127
$ws.Cells.Item(5, 18).Value = 0.04266859229641996
$ws.Cells.Item(5, 19).Value = 0.0418419654739287
$ws.Cells.Item(5, 20).Value = 0.0417672493930584
$ws.Cells.Item(5, 21).Value = 0.04141984305360232
$ws.Cells.Item(5, 22).Value = 0.04119146530959913
$ws.Cells.Item(5, 23).Value = 0.04103490510534425
$ws.Cells.Item(5, 24).Value = 0.04084265064425651
$ws.Cells.Item(5, 25).Value = 0.04073154009879833
$ws.Cells.Item(6, 3).Value = 1.627998352050781
$ws.Cells.Item(6, 5).Value = 2054.130812007405
$ws.Cells.Item(6, 6).Value = 0.09397968727101226
$ws.Cells.Item(6, 7).Value = 0.07725029474226243
$ws.Cells.Item(6, 8).Value = 0.06765354038788966
$ws.Cells.Item(6, 9).Value = 0.05969421596489668
$ws.Cells.Item(6, 10).Value = 0.05667148833115187
$ws.Cells.Item(6, 11).Value = 0.05218857474319434
$ws.Cells.Item(6, 12).Value = 0.04991447901796406
$ws.Cells.Item(6, 13).Value = 0.04797931778591412
$ws.Cells.Item(6, 14).Value = 0.04579911145361162
$ws.Cells.Item(6, 15).Value = 0.04419843789485547
$ws.Cells.Item(6, 16).Value = 0.04341355318344636
$ws.Cells.Item(6, 17).Value = 0.0423574932766021
$ws.Cells.Item(6, 18).Value = 0.04227656135565593
$ws.Cells.Item(6, 19).Value = 0.04150544231425234
$ws.Cells.Item(6, 20).Value = 0.04098740404844013
$ws.Cells.Item(6, 21).Value = 0.04072483718646696
$ws.Cells.Item(6, 22).Value = 0.04043371844877375
$ws.Cells.Item(6, 23).Value = 0.04025961650266923
$ws.Cells.Item(6, 24).Value = 0.04015532425913857
$ws.Cells.Item(6, 25).Value = 0.04004153629644062
$ws.Cells.Item(7, 3).Value = 1.660645961761475
$ws.Cells.Item(7, 5).Value = 2080.804101792068
$ws.Cells.Item(7, 6).Value = 0.09462290128124194
$ws.Cells.Item(7, 7).Value = 0.07677525461014907
$ws.Cells.Item(7, 8).Value = 0.06641138510167439
$ws.Cells.Item(7, 9).Value = 0.05942722586457996
$ws.Cells.Item(7, 10).Value = 0.05515482836843017
$ws.Cells.Item(7, 11).Value = 0.05067868193078766
$ws.Cells.Item(7, 12).Value = 0.04867569241981785
$ws.Cells.Item(7, 13).Value = 0.04641493210461217
$ws.Cells.Item(7, 14).Value = 0.04533032982887575
$ws.Cells.Item(7, 15).Value = 0.04473128211421507
$ws.Cells.Item(7, 16).Value = 0.0434970099806388
$ws.Cells.Item(7, 17).Value = 0.04271831683150806
$ws.Cells.Item(7, 18).Value = 0.04213455066135843
$ws.Cells.Item(7, 19).Value = 0.04168572635456904
$ws.Cells.Item(7, 20).Value = 0.04127686798007848
$ws.Cells.Item(7, 21).Value = 0.04104290570028814
$ws.Cells.Item(7, 22).Value = 0.04097753397229002
$ws.Cells.Item(7, 23).Value = 0.04080035965961794
$ws.Cells.Item(7, 24).Value = 0.04064207490583113
$ws.Cells.Item(7, 25).Value = 0.0405614834657323
$ws.Cells.Item(8, 3).Value = 1.712998151779175
$ws.Cells.Item(8, 5).Value = 2045.925541434346
$ws.Cells.Item(8, 6).Value = 0.09437991918576265
$ws.Cells.Item(8, 7).Value = 0.07522363841588058
$ws.Cells.Item(8, 8).Value = 0.06903016987404968
$ws.Cells.Item(8, 9).Value = 0.05822951539529845
$ws.Cells.Item(8, 10).Value = 0.05452740994649019
$ws.Cells.Item(8, 11).Value = 0.05148241424973311
$ws.Cells.Item(8, 12).Value = 0.04922276264525289
$ws.Cells.Item(8, 13).Value = 0.0472945582444373
$ws.Cells.Item(8, 14).Value = 0.04543844009585785
$ws.Cells.Item(8, 15).Value = 0.04406579705106247
$ws.Cells.Item(8, 16).Value = 0.04322218954773902
$ws.Cells.Item(8, 17).Value = 0.04238540932716396
$ws.Cells.Item(8, 18).Value = 0.04165709345687459
$ws.Cells.Item(8, 19).Value = 0.04109895929616249
$ws.Cells.Item(8, 20).Value = 0.04048934042370875
$ws.Cells.Item(8, 21).Value = 0.04039106012279217
$ws.Cells.Item(8, 22).Value = 0.04024102749323697
$ws.Cells.Item(8, 23).Value = 0.04008816599833929
$ws.Cells.Item(8, 24).Value = 0.03991005027198264
$ws.Cells.Item(8, 25).Value = 0.03988158950164417
$ws.Cells.Item(9, 3).Value = 1.629999399185181
$ws.Cells.Item(9, 5).Value = 2095.226719795983
$ws.Cells.Item(9, 6).Value = 0.09030333871940428
$ws.Cells.Item(9, 7).Value = 0.07745575044779146
$ws.Cells.Item(9, 8).Value = 0.06793274200817161
$ws.Cells.Item(9, 9).Value = 0.06120037129551467
$ws.Cells.Item(9, 10).Value = 0.0561405727252252
$ws.Cells.Item(9, 11).Value = 0.05227842319425915
$ws.Cells.Item(9, 12).Value = 0.05039184592087559
$ws.Cells.Item(9, 13).Value = 0.04698084159316915
$ws.Cells.Item(9, 14).Value = 0.04606839452637098
$ws.Cells.Item(9, 15).Value = 0.04551310504781189
$ws.Cells.Item(9, 16).Value = 0.04460866761615089
$ws.Cells.Item(9, 17).Value = 0.04358081994846193
$ws.Cells.Item(9, 18).Value = 0.04313716816574954
$ws.Cells.Item(9, 19).Value = 0.0421695602064293
$ws.Cells.Item(9, 20).Value = 0.04184168845057011
$ws.Cells.Item(9, 21).Value = 0.04156479875595883
$ws.Cells.Item(9, 22).Value = 0.04127823897516132
$ws.Cells.Item(9, 23).Value = 0.04095173095237709
$ws.Cells.Item(9, 24).Value = 0.04094487252950358
$ws.Cells.Item(9, 25).Value = 0.04084262611688074
$ws.Cells.Item(10, 3).Value = 1.483002662658691
$ws.Cells.Item(10, 5).Value = 2097.831214728378
$ws.Cells.Item(10, 6).Value = 0.09302231987941975
$ws.Cells.Item(10, 7).Value = 0.07409725264055522
$ws.Cells.Item(10, 8).Value = 0.06455503120587008
$ws.Cells.Item(10, 9).Value = 0.05837973567001255
$ws.Cells.Item(10, 10).Value = 0.055195758021547
$ws.Cells.Item(10, 11).Value = 0.0529182560631328
$ws.Cells.Item(10, 12).Value = 0.04921150677898703
$ws.Cells.Item(10, 13).Value = 0.04800088064283262
$ws.Cells.Item(10, 14).Value = 0.04742214413890924
$ws.Cells.Item(10, 15).Value = 0.04640997321981414
$ws.Cells.Item(10, 16).Value = 0.04510979918443056
$ws.Cells.Item(10, 17).Value = 0.04439140004018894
$ws.Cells.Item(10, 18).Value = 0.04300890919246345
$ws.Cells.Item(10, 19).Value = 0.04247786915420945
$ws.Cells.Item(10, 20).Value = 0.04196897224612891
$ws.Cells.Item(10, 21).Value = 0.04140335800212637
$ws.Cells.Item(10, 22).Value = 0.04128990609187624
$ws.Cells.Item(10, 23).Value = 0.04113509050208477
$ws.Cells.Item(10, 24).Value = 0.04097565638303467
$ws.Cells.Item(10, 25).Value = 0.04089339599860384
$ws.Cells.Item(11, 3).Value = 1.701997518539429
$ws.Cells.Item(11, 5).Value = 2078.782217955319
$ws.Cells.Item(11, 6).Value = 0.0951994750013921
$ws.Cells.Item(11, 7).Value = 0.07713783526644546
$ws.Cells.Item(11, 8).Value = 0.06699080052867359
$ws.Cells.Item(11, 9).Value = 0.06096992798555035
$ws.Cells.Item(11, 10).Value = 0.05530454330697981
$ws.Cells.Item(11, 11).Value = 0.05211003979437823
$ws.Cells.Item(11, 12).Value = 0.0483631760078902
$ws.Cells.Item(11, 13).Value = 0.04750738401897126
$ws.Cells.Item(11, 14).Value = 0.04551716737972711
$ws.Cells.Item(11, 15).Value = 0.04378971934256529
$ws.Cells.Item(11, 16).Value = 0.04320516558474126
$ws.Cells.Item(11, 17).Value = 0.04266311580957934
$ws.Cells.Item(11, 18).Value = 0.04208199836310199
$ws.Cells.Item(11, 19).Value = 0.04161073614020045
$ws.Cells.Item(11, 20).Value = 0.04145980797690887
$ws.Cells.Item(11, 21).Value = 0.04121624059241891
$ws.Cells.Item(11, 22).Value = 0.04095652452183285
$ws.Cells.Item(11, 23).Value = 0.04065991185997406
$ws.Cells.Item(11, 24).Value = 0.04054490363848002
$ws.Cells.Item(11, 25).Value = 0.0405220705254448
